# Insert a "," run right after the existing "sección" run in the second
# paragraph ("En esta sección se maquetará los datos de la factura. ...").
#
# A plain Find/Replace or Range.InsertAfter on this runtime coalesces every
# run in the paragraph that shares identical run-formatting (rPr) into one
# big run as part of any text-mutating edit. To keep the document's other,
# untouched runs split exactly as they were (matching the target diff) we
# "pin" every existing run in the paragraph with a harmless, transient
# Bold toggle (On then immediately back Off does NOT trigger the coalescing
# pass) so that, at the moment the new comma is spliced in, no two adjacent
# runs share identical formatting and none of them get merged. Afterwards
# we revert the temporary Bold flags (also a no-merge operation) and the
# paragraph ends up with exactly one new run - the comma - while every
# original run boundary is preserved.

$d = $word.ActiveDocument

# Anchor on the unique phrase that starts the target paragraph.
$anchor = $d.Content
$anchor.Find.Execute("En esta sección", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$pstart = $anchor.Start

# Lengths (in characters) of the paragraph's original runs, in order, taken
# from the document's current OOXML:
#   "En esta " / "sección" / " se maquetará " / "los datos de la factura. " /
#   "Como ya se " / "había " / "mencionado antes," / " el archivo " /
#   "getInvoice.js" / ", se utiliza en el componente " / "InvoiceApp" /
#   " como una función de servicio" / ", el cual se encarga de retornar" /
#   " la información para mostrarla en el documento HTML. "
$lens = @(8, 7, 14, 25, 11, 6, 17, 12, 13, 30, 10, 29, 32, 53)

$starts = @()
$pos = $pstart
foreach ($len in $lens) {
    $starts += $pos
    $pos += $len
}

# Pin every other run (odd index) as Bold=On so that no two adjacent runs
# ever share identical formatting while we splice the new text in.
for ($i = 0; $i -lt $lens.Length; $i++) {
    if ($i % 2 -eq 1) {
        $s = $starts[$i]
        $e = $s + $lens[$i]
        $pin = $d.Range($s, $e)
        $pin.Bold = $true
    }
}

# Insert the new "," immediately after "sección" (run index 1), i.e. right
# before the run that currently reads " se maquetará ".
$insPos = $starts[2]
$ins = $d.Range($insPos, $insPos)
$ins.InsertAfter(",")

# Everything from the insertion point onward shifted right by one character.
for ($i = 2; $i -lt $lens.Length; $i++) {
    $starts[$i] = $starts[$i] + 1
}

# Revert the temporary Bold pin on the odd-indexed runs.
for ($i = 0; $i -lt $lens.Length; $i++) {
    if ($i % 2 -eq 1) {
        $s = $starts[$i]
        $e = $s + $lens[$i]
        $pin = $d.Range($s, $e)
        $pin.Bold = $false
    }
}

# The new comma inherited Bold from its (pinned) left neighbour; clear it.
$comma = $d.Range($insPos, $insPos + 1)
$comma.Bold = $false
